# The workbook contained a spurious row (row 6) holding only the header
# label "grandes regiões e unidades da federação" with no data underneath
# it (the actual "norte" data row that follows already carries the real
# regional figures). Remove that stray row entirely: this deletes its
# label from the shared strings table, shifts every subsequent row up by
# one, and drops the last (now-empty) row from the bottom of the range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Delete()
